# Add "Wellcome Trust, " to the list of funders whose grants the author
# reviews, right before "UK Research and Innovation: ".
$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute(
    "UK Research and Innovation: ",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "",
    0
)

if ($found) {
    $rng.InsertBefore("Wellcome Trust, ")
}

Write-Output "Found and inserted: $found"
